# isotope model version update
# - replace the "ND" placeholder text entries in F4 and B7 with their
#   now-available measured numeric values
# - move the active selection to B8 (last cell touched while reviewing data)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F4").Value = -2.5
$ws.Range("B7").Value = 3

$ws.Range("B8").Select() | Out-Null
